# Trade #6 closed at 2026-02-17 07:52:46 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.05   # Current Capital
$wsSummary.Range("B4").Value = 0.05      # Total P&L $
$wsSummary.Range("B5").Value = 0.17      # Total P&L %
$wsSummary.Range("B6").Value = 6         # Total Trades
$wsSummary.Range("B7").Value = 3         # Winning Trades
$wsSummary.Range("B9").Value = 50        # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.05     # Capital
$wsStatus.Range("D4").Value = 6          # Trades
$wsStatus.Range("E4").Value = 0.05       # P&L $
$wsStatus.Range("F4").Value = 0.05       # P&L %
$wsStatus.Range("G4").Value = 50         # Win Rate %

# --- New trade row (#6) shared by "All Trades" and "MarketMaking" sheets ---
$tradeValues = @{
    A = 6
    B = "2026-02-17"
    C = "07:52:40"
    D = "MarketMaking"
    E = "UP"
    F = 0.6
    G = 0.73
    H = "CLOSED"
    I = 21.6667
    J = 0.13
    K = 100.05
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($col in $tradeValues.Keys) {
        $cell = $ws.Range("$col`7")
        $value = $tradeValues[$col]
        if ($col -eq "B") {
            # "2026-02-17" looks like a date to Excel's auto-detection; force
            # it to be stored as literal text (matching the other rows) by
            # entering it with a leading apostrophe and then clearing the
            # resulting "quote prefix" formatting back to the default style.
            $cell.Value = "'$value"
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
